$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.148.85"
$ws.Range("E2").Value = "  +0.89%  "

$ws.Range("D3").Value = "2.928.82"
$ws.Range("E3").Value = "  +0.99%  "

$ws.Range("E4").Value = "  -0.15%  "

$ws.Range("D5").Value = "'591.80"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").Value = "'145.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.22%  "

$ws.Range("E7").Value = "  -0.01%  "

$ws.Range("E8").Value = "  +0.74%  "

$ws.Range("D9").Value = "'6.93"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +4.01%  "

$ws.Range("E10").Value = "  +0.23%  "

$ws.Range("E11").Value = "  -0.25%  "

$ws.Range("D12").Value = "'0.0000226"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.58%  "

$ws.Range("D13").Value = "'33.73"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.29%  "

$ws.Range("E14").Value = "  -0.27%  "

$ws.Range("D15").Value = "3.412.60"
$ws.Range("E15").Value = "  +0.76%  "

$ws.Range("D16").Value = "60.986.01"
$ws.Range("E16").Value = "  +0.53%  "

$ws.Range("E17").Value = "  +0.96%  "

$ws.Range("D18").Value = "2.929.88"
$ws.Range("E18").Value = "  +0.84%  "

$ws.Range("D19").Value = "'437.14"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.42%  "

$ws.Range("D20").Value = "'13.45"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("E21").Value = "  -0.32%  "

$ws.Range("D22").Value = "'7.12"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.13%  "

$ws.Range("D23").Value = "'81.57"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.04%  "

$ws.Range("D24").Value = "'11.08"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.84%  "

$ws.Range("D25").Value = "'2.21"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.16%  "

$ws.Range("D26").Value = "'11.89"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.43%  "

$ws.Range("E27").Value = "  -0.03%  "

$ws.Range("D28").Value = "'2.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.38%  "

$ws.Range("E29").Value = "  +0.12%  "

$ws.Range("E30").Value = "  -1.85%  "

$ws.Range("D31").Value = "'0.111"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.59%  "

$ws.Range("D32").Value = "'26.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.31%  "

$ws.Range("E33").Value = "  -0.02%  "

$ws.Range("D34").Value = "0.0₃0869"
$ws.Range("E34").Value = "  +2.58%  "

$ws.Range("E35").Value = "  +0.55%  "

$ws.Range("D36").Value = "'5.65"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.47%  "

$ws.Range("D37").Value = "'3.01"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.13%  "

$ws.Range("E40").Value = "  +0.58%  "

$ws.Range("D41").Value = "'42.12"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.70%  "

$ws.Range("D42").Value = "'0.287"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "

$ws.Range("D43").Value = "'375.06"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.97%  "

$ws.Range("E44").Value = "  -0.43%  "

$ws.Range("D45").Value = "2.691.00"
$ws.Range("E45").Value = "  -0.15%  "

$ws.Range("D46").Value = "'133.46"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.09%  "

$ws.Range("E47").Value = "  -0.05%  "

$ws.Range("E48").Value = "  -0.15%  "

$ws.Range("E49").Value = "  -0.08%  "

$ws.Range("E50").Value = "  -1.54%  "

$ws.Range("E51").Value = "  +1.19%  "

# Row 38: Kaspa -> Stacks (rows swap ranking position)
$ws.Range("B38").Value = "Stacks"
$ws.Range("C38").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D38").Value = "'1.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.00%  "

# Row 39: Stacks -> Kaspa
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "'0.123"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.18%  "
